# Insert a new weekly price row at row 551 ("1a nueva(o)" quality reading
# for 2023-12-07), pushing the existing rows 551:579 down to 552:580.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(551).Insert()

$ws.Cells.Item(551, 1).Value = 4
$ws.Cells.Item(551, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(551, 3).Value = "Los Lagos"
$ws.Cells.Item(551, 4).Value = 45267
$ws.Cells.Item(551, 5).Value = 10
$ws.Cells.Item(551, 6).Value = 100112045
$ws.Cells.Item(551, 7).Value = "Zapallo"
$ws.Cells.Item(551, 8).Value = "Paine"
$ws.Cells.Item(551, 9).Value = "1a nueva(o)"
$ws.Cells.Item(551, 10).Value = 500
$ws.Cells.Item(551, 11).Value = 1500
$ws.Cells.Item(551, 12).Value = 1500
$ws.Cells.Item(551, 13).Value = 1500
$ws.Cells.Item(551, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(551, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(551, 16).Value = 1500
$ws.Cells.Item(551, 17).Value = 1
$ws.Cells.Item(551, 18).Value = "Hortaliza"
